$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 0.03299100410245547
$ws.Range("C3").Value = 0.03509343290161347
$ws.Range("D3").Value = 0.03351395295910102

$ws.Range("C4").Value = 0.033372451942018

$ws.Range("B5").Value = 0.05856157327945263
$ws.Range("C5").Value = 0.05791343948902639
$ws.Range("D5").Value = 0.06176581628842315
